# Applies scheduled market-data refresh values to the profit-tracking sheets.
$wb = $excel.ActiveWorkbook

# ---- Worksheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1164.7333
$ws.Range("I19").Value = 1664.5625
$ws.Range("J19").Value = 593.5
$ws.Range("K19").Value = 1664.5625
$ws.Range("L19").Value = 593.5
$ws.Range("M19").Value = -1489.5625
$ws.Range("N19").Value = -943.5
$ws.Range("H62").Value = 2545.2563
$ws.Range("I62").Value = 1987.5
$ws.Range("J62").Value = 2689.1936
$ws.Range("K62").Value = 1987.5
$ws.Range("L62").Value = 2689.1936
$ws.Range("M62").Value = -1363.5
$ws.Range("N62").Value = -3937.1936
$ws.Range("H65").Value = 2545.2563
$ws.Range("I65").Value = 1987.5
$ws.Range("J65").Value = 2689.1936
$ws.Range("K65").Value = 9937.5
$ws.Range("L65").Value = 13445.968
$ws.Range("M65").Value = -6817.5
$ws.Range("N65").Value = -19685.968
$ws.Range("H113").Value = 2424
$ws.Range("I113").Value = 2448.3333
$ws.Range("J113").Value = 2399.6667
$ws.Range("K113").Value = 2448.3333
$ws.Range("L113").Value = 2399.6667
$ws.Range("M113").Value = 805.6667000000002
$ws.Range("N113").Value = -8907.6667
$ws.Range("H127").Value = 1133.1333
$ws.Range("I127").Value = 799.6667
$ws.Range("J127").Value = 1633.3334
$ws.Range("K127").Value = 2399.0001
$ws.Range("L127").Value = 4900.0002
$ws.Range("M127").Value = 2560.9999
$ws.Range("N127").Value = -14820.0002
$ws.Range("H132").Value = 336651.6
$ws.Range("I132").Value = 336651.6
$ws.Range("K132").Value = 1009954.8
$ws.Range("M132").Value = -1007424.8

# ---- Worksheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 977.9286
$ws.Range("I45").Value = 958.7143
$ws.Range("J45").Value = 997.1429000000001
$ws.Range("K45").Value = 958.7143
$ws.Range("L45").Value = 997.1429000000001
$ws.Range("M45").Value = -581.7143
$ws.Range("N45").Value = -1751.1429

# ---- Worksheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1744.4
$ws.Range("I20").Value = 1774.381
$ws.Range("J20").Value = 1674.4445
$ws.Range("K20").Value = 1774.381
$ws.Range("L20").Value = 1674.4445
$ws.Range("M20").Value = -1527.381
$ws.Range("N20").Value = -2168.4445

# ---- Worksheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1843.4584
$ws.Range("I58").Value = 1529.6875
$ws.Range("J58").Value = 2471
$ws.Range("K58").Value = 1529.6875
$ws.Range("L58").Value = 2471
$ws.Range("M58").Value = -1326.6875
$ws.Range("N58").Value = -2877
$ws.Range("H62").Value = 2308.1538
$ws.Range("I62").Value = 2180
$ws.Range("J62").Value = 2735.3333
$ws.Range("K62").Value = 2180
$ws.Range("L62").Value = 2735.3333
$ws.Range("M62").Value = -1556
$ws.Range("N62").Value = -3983.3333
$ws.Range("H65").Value = 2308.1538
$ws.Range("I65").Value = 2180
$ws.Range("J65").Value = 2735.3333
$ws.Range("K65").Value = 10900
$ws.Range("L65").Value = 13676.6665
$ws.Range("M65").Value = -7780
$ws.Range("N65").Value = -19916.6665
$ws.Range("H86").Value = 2113.5715
$ws.Range("I86").Value = 1871.125
$ws.Range("J86").Value = 2262.7693
$ws.Range("K86").Value = 1871.125
$ws.Range("L86").Value = 2262.7693
$ws.Range("M86").Value = -748.125
$ws.Range("N86").Value = -4508.7693
$ws.Range("H89").Value = 2113.5715
$ws.Range("I89").Value = 1871.125
$ws.Range("J89").Value = 2262.7693
$ws.Range("K89").Value = 9355.625
$ws.Range("L89").Value = 11313.8465
$ws.Range("M89").Value = -3739.625
$ws.Range("N89").Value = -22545.8465
$ws.Range("H99").Value = 1506
$ws.Range("I99").Value = 1506
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1506
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 2234.44
$ws.Range("I122").Value = 2147.3845
$ws.Range("J122").Value = 2328.75
$ws.Range("K122").Value = 6442.1535
$ws.Range("L122").Value = 6986.25
$ws.Range("M122").Value = -3992.1535
$ws.Range("N122").Value = -11886.25
$ws.Range("H126").Value = 1506
$ws.Range("I126").Value = 1506
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4518
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2048
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2197.0312
$ws.Range("I132").Value = 1590.05
$ws.Range("K132").Value = 4770.15
$ws.Range("M132").Value = -2240.15
$ws.Range("H136").Value = 1843.4584
$ws.Range("I136").Value = 1529.6875
$ws.Range("J136").Value = 2471
$ws.Range("K136").Value = 4589.0625
$ws.Range("L136").Value = 7413
$ws.Range("M136").Value = -2039.0625
$ws.Range("N136").Value = -12513

# ---- Worksheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 534212.5
$ws.Range("I2").Value = 42.444443
$ws.Range("J2").Value = 854714.5600000001
$ws.Range("K2").Value = 254.666658
$ws.Range("L2").Value = 5128287.36
$ws.Range("M2").Value = -141.666658
$ws.Range("N2").Value = -5128513.36
$ws.Range("H98").Value = 537.45
$ws.Range("J98").Value = 460
$ws.Range("L98").Value = 1380
$ws.Range("N98").Value = -4376
$ws.Range("H101").Value = 6850
$ws.Range("J101").Value = 6850
$ws.Range("L101").Value = 20550
$ws.Range("N101").Value = -25418
$ws.Range("H131").Value = 1113043.1
$ws.Range("J131").Value = 1267342.5
$ws.Range("L131").Value = 3802027.5
$ws.Range("N131").Value = -3812107.5

# ---- Worksheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064
$ws.Range("H132").Value = 2616.4856
$ws.Range("I132").Value = 2458.4348
$ws.Range("J132").Value = 2919.4167
$ws.Range("K132").Value = 7375.3044
$ws.Range("L132").Value = 8758.250100000001
$ws.Range("M132").Value = -4845.3044
$ws.Range("N132").Value = -13818.2501

# ---- Worksheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1457.4
$ws.Range("I7").Value = 1311.4667
$ws.Range("J7").Value = 1603.3334
$ws.Range("K7").Value = 1311.4667
$ws.Range("L7").Value = 1603.3334
$ws.Range("M7").Value = -1199.4667
$ws.Range("N7").Value = -1827.3334
$ws.Range("H16").Value = 1642.96
$ws.Range("I16").Value = 1682.25
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1682.25
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -1512.25
$ws.Range("N16").Value = -1040
$ws.Range("H40").Value = 1508.04
$ws.Range("I40").Value = 1299.7778
$ws.Range("J40").Value = 2043.5714
$ws.Range("K40").Value = 1299.7778
$ws.Range("L40").Value = 2043.5714
$ws.Range("M40").Value = -1163.7778
$ws.Range("N40").Value = -2315.5714
$ws.Range("H126").Value = 1457.4
$ws.Range("I126").Value = 1311.4667
$ws.Range("J126").Value = 1603.3334
$ws.Range("K126").Value = 3934.4001
$ws.Range("L126").Value = 4810.0002
$ws.Range("M126").Value = -1464.4001
$ws.Range("N126").Value = -9750.0002

# ---- Worksheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10049.143
$ws.Range("I62").Value = 9416.333000000001
$ws.Range("J62").Value = 10523.75
$ws.Range("K62").Value = 9416.333000000001
$ws.Range("L62").Value = 10523.75
$ws.Range("M62").Value = -8792.333000000001
$ws.Range("N62").Value = -11771.75
$ws.Range("H65").Value = 10049.143
$ws.Range("I65").Value = 9416.333000000001
$ws.Range("J65").Value = 10523.75
$ws.Range("K65").Value = 47081.665
$ws.Range("L65").Value = 52618.75
$ws.Range("M65").Value = -43961.665
$ws.Range("N65").Value = -58858.75
$ws.Range("H113").Value = 407.09525
$ws.Range("I113").Value = 448.08334
$ws.Range("K113").Value = 1344.25002
$ws.Range("M113").Value = 825.7499800000001
$ws.Range("H132").Value = 1249.0222
$ws.Range("I132").Value = 675.7273
$ws.Range("J132").Value = 2825.5833
$ws.Range("K132").Value = 2027.1819
$ws.Range("L132").Value = 8476.749899999999
$ws.Range("M132").Value = 502.8181
$ws.Range("N132").Value = -13536.7499
$ws.Range("H136").Value = 7737.8066
$ws.Range("I136").Value = 8482.27
$ws.Range("J136").Value = 3866.6
$ws.Range("K136").Value = 25446.81
$ws.Range("L136").Value = 11599.8
$ws.Range("M136").Value = -22896.81
$ws.Range("N136").Value = -16699.8

Write-Host "Applied 210 cell updates across 8 worksheets."
